# The workbook was opened, reviewed ("tested"), and re-saved without any
# data changes. The only user-visible effect captured in the diff is that
# the active selection on Sheet1 ends up at the default cell A1 instead of
# the previously-saved C5.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()
$ws.Range("A1").Select()
